$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 815 (shifts existing rows 815:856 down to 816:857,
# and bumps the sheet's used-range dimension from D856 to D857).
$ws.Rows.Item(815).Insert()

# The new row 815 holds an extra data point for 2026/02/19 (Thursday),
# time 14, ranking 79 - matching the same text formatting used by the
# other date/weekday cells in the column (plain text, no number format).
$ws.Cells.Item(815, 1).NumberFormat = "@"
$ws.Cells.Item(815, 1).Value = "2026/02/19"
$ws.Cells.Item(815, 1).ClearFormats()
$ws.Cells.Item(815, 2).Value = "木"
$ws.Cells.Item(815, 3).Value = 14
$ws.Cells.Item(815, 4).Value = 79
